# "add 2 examples on Roba che non funziona"
#
# 1) Remove the stray "_GoBack" bookmark that sat at the end of the
#    "Prenotazione (Conferma) -> risolto" paragraph.
# 2) Append two new "Motivo Prenotazione (...)" list items to the
#    "ESEMPI AGGIUNTI" section (right after the existing
#    "Prenotazione (recuperaCamerePerDate)" item), followed by a new
#    (otherwise empty) list paragraph that now carries the "_GoBack"
#    bookmark - mirroring where Word's cursor/last-edit marker ends up
#    after typing the new lines.

$d = $word.ActiveDocument

# --- 1. drop the old _GoBack bookmark -------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. insert the two new examples + the paragraph that now holds the ----
#        _GoBack bookmark, right before the trailing empty paragraph
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertPoint = $lastPara.Range
$insertPoint.Collapse(1)  # wdCollapseStart

$newParasXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Motivo Prenotazione (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>getAllMotiviPrenotazione</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Motivo Prenotazione (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>getMotivoPrenotazioneById</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$insertPoint.InsertXML($newParasXml)

Write-Output "done"
